$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Acelga" (Chard) series at
# Macroferia Regional de Talca. It belongs right before the current row 149
# (it is the most recent date among the existing entries), so insert a new
# row there; every row from the old 149 through 180 shifts down by one
# (to 150-181), keeping all of their original values intact.
$ws.Rows(149).Insert()

# Populate the freshly inserted row 149 with the new data point. It mirrors
# the (now shifted) row 150 in every column except the date (D), which is
# the new reporting date.
$ws.Range("A149").Value = 5
$ws.Range("B149").Value = "Macroferia Regional de Talca"
$ws.Range("C149").Value = "Maule"
$ws.Range("D149").Value = 44476
$ws.Range("E149").Value = 7
$ws.Range("F149").Value = 100112009
$ws.Range("G149").Value = "Acelga"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 500
$ws.Range("K149").Value = 2000
$ws.Range("L149").Value = 2000
$ws.Range("M149").Value = 2000
$ws.Range("N149").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O149").Value = "Región del Maule"
$ws.Range("P149").Value = 500
$ws.Range("Q149").Value = 4
$ws.Range("R149").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the column
# (same number format as D150, etc.).
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat
